# issue #5: stock data from json to db
# Add a "category" column (value "normal") and "source_file" / "index"
# columns to the 股票 (stock) worksheet, matching the new json->db export
# format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (column I).
# This shifts date/legislator_name/legislator_id from I/J/K to J/K/L.
$ws.Columns.Item(9).Insert()

# --- New column I: category -------------------------------------------------
$ws.Range("I1").Value = "category"

# --- New column M: source_file ----------------------------------------------
$ws.Range("M1").Value = "source_file"

# --- New column N: index -----------------------------------------------------
$ws.Range("N1").Value = "index"

# --- Data rows ---------------------------------------------------------------
$rows = @(2, 3, 4, 5, 6)
$indexValues = @{2 = 68; 3 = 69; 4 = 70; 5 = 72; 6 = 73}

foreach ($r in $rows) {
    $ws.Range("I$r").Value = "normal"
}

foreach ($r in $rows) {
    $ws.Range("M$r").Value = "tmp82d01"
}

foreach ($r in $rows) {
    $ws.Range("N$r").Value = $indexValues[$r]
}
